$wb = $excel.ActiveWorkbook

# Rename the strain sheets to include the "_log2_expression" suffix
$wb.Worksheets.Item("wt").Name = "wt_log2_expression"
$wb.Worksheets.Item("dcin5").Name = "dcin5_log2_expression"

# Make the dcin5_log2_expression sheet the active / selected sheet
$wb.Worksheets.Item("dcin5_log2_expression").Activate()
